# Apply updated crypto price/volume figures and re-ordered coin rows
# per the "Updated symbol list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.31"
$ws.Range("E2").Value = "'0.82%"
$ws.Range("D3").Value = "'27.27"
$ws.Range("E3").Value = "'2.66%"
$ws.Range("D4").Value = "'4.884"
$ws.Range("E4").Value = "'0.32%"
$ws.Range("D5").Value = "'0.06416"
$ws.Range("E5").Value = "'1.19%"
$ws.Range("D6").Value = "'6.964"
$ws.Range("E6").Value = "'1.10%"
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = "'3.321"
$ws.Range("E7").Value = "'0.21%"
$ws.Range("D8").Value = "'1.177"
$ws.Range("E8").Value = "'-6.19%"
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = "'0.8840"
$ws.Range("E9").Value = "'1.89%"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1549"
$ws.Range("E10").Value = "'-0.31%"
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = "'0.05127"
$ws.Range("E11").Value = "'-1.72%"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.07442"
$ws.Range("E12").Value = "'0.39%"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.02900"
$ws.Range("E13").Value = "'-0.76%"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.08978"
$ws.Range("E14").Value = "'-0.57%"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = "'0.001563"
$ws.Range("E15").Value = "'-0.67%"
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = "'0.0006376"
$ws.Range("E16").Value = "'0.69%"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = "'0.006167"
$ws.Range("E17").Value = "'3.39%"
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = "'3.477"
$ws.Range("E18").Value = "'0.80%"
$ws.Range("D19").Value = "'2.274"
$ws.Range("E19").Value = "'0.07%"
$ws.Range("D21").Value = "'0.1331"
$ws.Range("E21").Value = "'-0.23%"
$ws.Range("D22").Value = "'3.903"
$ws.Range("E22").Value = "'-0.14%"
$ws.Range("D23").Value = "'0.04432"
$ws.Range("E23").Value = "'1.38%"
$ws.Range("E24").Value = "'8.76%"
$ws.Range("D26").Value = "'0.001178"
$ws.Range("E26").Value = "'0.13%"
$ws.Range("D27").Value = "'0.003865"
$ws.Range("E28").Value = "'-1.60%"
$ws.Range("E29").Value = "'15.66%"
$ws.Range("D40").Value = "'0.04149"
$ws.Range("E40").Value = "'1.01%"
$ws.Range("D41").Value = "'0.006789"
$ws.Range("E41").Value = "'-1.44%"
$ws.Range("E42").Value = "'0.68%"
$ws.Range("D43").Value = "'0.002001"
$ws.Range("E43").Value = "'-6.64%"
$ws.Range("D44").Value = "'0.01147"
$ws.Range("E44").Value = "'6.90%"
$ws.Range("D45").Value = "'0.00005315"
$ws.Range("E45").Value = "'1.20%"
$ws.Range("D46").Value = "'1.685"
$ws.Range("E46").Value = "'13.33%"
$ws.Range("D47").Value = "'0.01854"
$ws.Range("E47").Value = "'-7.27%"
